# Populate the "Results" sheet with the list of coffee shops in the
# Linh Đàm area (rows 2-10), below the existing header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the new data range to be stored as text so values that look
# numeric ("105.8251937", "4,4", "(165)", ...) keep their literal
# representation instead of being auto-converted to real numbers.
$ws.Range("A2:G10").NumberFormat = "@"

# Row 2
$ws.Range("A2").Value = "Linh Đàm"
$ws.Range("B2").Value = "Botanix Cafe"
$ws.Range("C2").Value = "105.8251937"
$ws.Range("D2").Value = "20.9656864"
$ws.Range("E2").Value = "cà phê"
$ws.Range("F2").Value = "4,4"
$ws.Range("G2").Value = "(165)"

# Row 3
$ws.Range("A3").Value = "Linh Đàm"
$ws.Range("B3").Value = "Cộng Cà Phê"
$ws.Range("C3").Value = "105.825746"
$ws.Range("D3").Value = "20.9634174"
$ws.Range("E3").Value = "cà phê"
$ws.Range("F3").Value = "4,0"
$ws.Range("G3").Value = "(1.281)"

# Row 4
$ws.Range("A4").Value = "Linh Đàm"
$ws.Range("B4").Value = "Cafe Chouchi"
$ws.Range("C4").Value = "105.8226406"
$ws.Range("D4").Value = "20.9633296"
$ws.Range("E4").Value = "cà phê"
$ws.Range("F4").Value = "4,2"
$ws.Range("G4").Value = "(443)"

# Row 5
$ws.Range("A5").Value = "Linh Đàm"
$ws.Range("B5").Value = "Quán The Coffee House"
$ws.Range("C5").Value = "105.8229694"
$ws.Range("D5").Value = "20.9634459"
$ws.Range("E5").Value = "cà phê"
$ws.Range("F5").Value = "4,4"
$ws.Range("G5").Value = "(535)"

# Row 6
$ws.Range("A6").Value = "Linh Đàm"
$ws.Range("B6").Value = "House of Cha Coffee"
$ws.Range("C6").Value = "105.8238937"
$ws.Range("D6").Value = "20.962895"
$ws.Range("E6").Value = "cà phê"
$ws.Range("F6").Value = "4,5"
$ws.Range("G6").Value = "(26)"

# Row 7
$ws.Range("A7").Value = "Linh Đàm"
$ws.Range("B7").Value = "Highlands Coffee"
$ws.Range("C7").Value = "105.8254817"
$ws.Range("D7").Value = "20.9638615"
$ws.Range("E7").Value = "cà phê"
$ws.Range("F7").Value = "4,2"
$ws.Range("G7").Value = "(867)"

# Row 8
$ws.Range("A8").Value = "Linh Đàm"
$ws.Range("B8").Value = "Laika Cafe Linh Đàm"
$ws.Range("C8").Value = "105.8256359"
$ws.Range("D8").Value = "20.9653021"
$ws.Range("E8").Value = "cà phê"
$ws.Range("F8").Value = "4,2"
$ws.Range("G8").Value = "(239)"

# Row 9
$ws.Range("A9").Value = "Linh Đàm"
$ws.Range("B9").Value = "Forli Coffee & Tea"
$ws.Range("C9").Value = "105.8272274"
$ws.Range("D9").Value = "20.9712171"
$ws.Range("E9").Value = "cà phê"
$ws.Range("F9").Value = "4,7"
$ws.Range("G9").Value = "(50)"

# Row 10
$ws.Range("A10").Value = "Linh Đàm"
$ws.Range("B10").Value = "Milano Coffee"
$ws.Range("C10").Value = "105.8251697"
$ws.Range("D10").Value = "20.966328"
$ws.Range("E10").Value = "cà phê"
$ws.Range("F10").Value = "4,2"
$ws.Range("G10").Value = "(132)"
